$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-DateText($cellAddr, $text) {
    # Force the cell to stay text (Excel otherwise auto-converts ambiguous
    # dd-mm-yyyy strings, where day <= 12, into real date serials). After
    # writing the text, restore the cell's style to Normal so no stray
    # number-format / style index is left behind on the cell.
    $ws.Range($cellAddr).NumberFormat = "@"
    $ws.Range($cellAddr).Value = $text
    $ws.Range($cellAddr).Style = "Normal"
}

# Update date labels (slash -> dash format) and attendance counters
# Columns: A=Date, D=Total Attendance Count, E=Real, F=Duplicate, G=Invalid, H=Absent

Set-DateText "A3" "28-07-2022"
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 1

Set-DateText "A4" "01-08-2022"
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("H4").Value = 0

Set-DateText "A5" "04-08-2022"
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1
$ws.Range("H5").Value = 0

Set-DateText "A6" "08-08-2022"
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 1
$ws.Range("H6").Value = 0

Set-DateText "A7" "11-08-2022"

Set-DateText "A8" "15-08-2022"

Set-DateText "A9" "18-08-2022"

Set-DateText "A10" "22-08-2022"

Set-DateText "A11" "25-08-2022"
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 1
$ws.Range("H11").Value = 0

Set-DateText "A12" "29-08-2022"
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 1
$ws.Range("H12").Value = 0

Set-DateText "A13" "01-09-2022"
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 1
$ws.Range("H13").Value = 0

Set-DateText "A14" "05-09-2022"

Set-DateText "A15" "08-09-2022"

Set-DateText "A16" "12-09-2022"

Set-DateText "A17" "15-09-2022"

Set-DateText "A18" "19-09-2022"

Set-DateText "A19" "22-09-2022"

Set-DateText "A20" "26-09-2022"

Set-DateText "A21" "29-09-2022"
